$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 20.97917533333333
$ws.Range("H2").Value = 62.93752600000001
$ws.Range("I2").Value = 0.2451892257562263
$ws.Range("J2").Value = 0.2451892257562263
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.521994666666667
$ws.Range("N2").Value = 7.565983999999999
$ws.Range("O2").Value = 0.01218715015226367
$ws.Range("P2").Value = 0.01218715015226367
$ws.Range("Q2").Value = 52.90936830173155
$ws.Range("R2").Value = 476.184314715584
$ws.Range("S2").Value = 0.002988157910008404
$ws.Range("T2").Value = 0.002988157910008404

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 20.97917533333333
$ws.Range("H3").Value = 62.93752600000001
$ws.Range("I3").Value = 0.2451892257562263
$ws.Range("J3").Value = 0.2451892257562263
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 157.1645253333333
$ws.Range("N3").Value = 471.493576
$ws.Range("O3").Value = 0.7594733225102963
$ws.Range("P3").Value = 0.7594733225102964
$ws.Range("Q3").Value = 3297.182133148108
$ws.Range("R3").Value = 29674.63919833298
$ws.Range("S3").Value = 0.1862146759288083
$ws.Range("T3").Value = 0.1862146759288083

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 20.97917533333333
$ws.Range("H4").Value = 62.93752600000001
$ws.Range("I4").Value = 0.2451892257562263
$ws.Range("J4").Value = 0.2451892257562263
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 47.252316
$ws.Range("N4").Value = 141.756948
$ws.Range("O4").Value = 0.2283395273374399
$ws.Range("P4").Value = 0.2283395273374399
$ws.Range("Q4").Value = 991.3146222700721
$ws.Range("R4").Value = 8921.831600430649
$ws.Range("S4").Value = 0.05598639191740955
$ws.Range("T4").Value = 0.05598639191740955

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 59.05285266666667
$ws.Range("H5").Value = 177.158558
$ws.Range("I5").Value = 0.6901664624076501
$ws.Range("J5").Value = 0.6901664624076501
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.521994666666667
$ws.Range("N5").Value = 7.565983999999999
$ws.Range("O5").Value = 0.01218715015226367
$ws.Range("P5").Value = 0.01218715015226367
$ws.Range("Q5").Value = 148.9309794767858
$ws.Range("R5").Value = 1340.378815291072
$ws.Range("S5").Value = 0.008411162307418672
$ws.Range("T5").Value = 0.008411162307418672

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 59.05285266666667
$ws.Range("H6").Value = 177.158558
$ws.Range("I6").Value = 0.6901664624076501
$ws.Range("J6").Value = 0.6901664624076501
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 157.1645253333333
$ws.Range("N6").Value = 471.493576
$ws.Range("O6").Value = 0.7594733225102963
$ws.Range("P6").Value = 0.7594733225102964
$ws.Range("Q6").Value = 9281.013558935934
$ws.Range("R6").Value = 83529.12203042342
$ws.Range("S6").Value = 0.5241630162899156
$ws.Range("T6").Value = 0.5241630162899157

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 59.05285266666667
$ws.Range("H7").Value = 177.158558
$ws.Range("I7").Value = 0.6901664624076501
$ws.Range("J7").Value = 0.6901664624076501
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 47.252316
$ws.Range("N7").Value = 141.756948
$ws.Range("O7").Value = 0.2283395273374399
$ws.Range("P7").Value = 0.2283395273374399
$ws.Range("Q7").Value = 2790.384054906776
$ws.Range("R7").Value = 25113.45649416098
$ws.Range("S7").Value = 0.1575922838103158
$ws.Range("T7").Value = 0.1575922838103158

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5.531174333333333
$ws.Range("H8").Value = 16.593523
$ws.Range("I8").Value = 0.06464431183612354
$ws.Range("J8").Value = 0.06464431183612354
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.521994666666667
$ws.Range("N8").Value = 7.565983999999999
$ws.Range("O8").Value = 0.01218715015226367
$ws.Range("P8").Value = 0.01218715015226367
$ws.Range("Q8").Value = 13.94959216907022
$ws.Range("R8").Value = 125.546329521632
$ws.Range("S8").Value = 0.0007878299348365931
$ws.Range("T8").Value = 0.0007878299348365931

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5.531174333333333
$ws.Range("H9").Value = 16.593523
$ws.Range("I9").Value = 0.06464431183612354
$ws.Range("J9").Value = 0.06464431183612354
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 157.1645253333333
$ws.Range("N9").Value = 471.493576
$ws.Range("O9").Value = 0.7594733225102963
$ws.Range("P9").Value = 0.7594733225102964
$ws.Range("Q9").Value = 869.3043886342498
$ws.Range("R9").Value = 7823.739497708249
$ws.Range("S9").Value = 0.04909563029157242
$ws.Range("T9").Value = 0.04909563029157243

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 5.531174333333333
$ws.Range("H10").Value = 16.593523
$ws.Range("I10").Value = 0.06464431183612354
$ws.Range("J10").Value = 0.06464431183612354
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 47.252316
$ws.Range("N10").Value = 141.756948
$ws.Range("O10").Value = 0.2283395273374399
$ws.Range("P10").Value = 0.2283395273374399
$ws.Range("Q10").Value = 261.360797449756
$ws.Range("R10").Value = 2352.247177047804
$ws.Range("S10").Value = 0.01476085160971452
$ws.Range("T10").Value = 0.01476085160971452
